# Blind model names in the Evaluations sheet, column C (rows 2-49)
# Mapping: claude-opus-4.5 -> Model A, gemini-3-pro -> Model B,
#          gpt-5.1 -> Model C, kimi-k2 -> Model D

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluations")

$map = @{
    "gpt-5.1"         = "Model C"
    "claude-opus-4.5" = "Model A"
    "gemini-3-pro"    = "Model B"
    "kimi-k2"         = "Model D"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    if ($null -ne $current -and $map.ContainsKey([string]$current)) {
        $cell.Value = $map[[string]$current]
    }
}

$wb.Save()
